{"js": "// Apply the dated-worksheet update: refresh the date line and every\n// two-digit multiplication problem in the table with the new values\n// from the commit.\nconst replacements = [\n  [\"2025-04-29 Tuesday\", \"2025-04-30 Wednesday\"],\n  [\"33\u00d774=\", \"49\u00d751=\"],\n  [\"24\u00d742=\", \"60\u00d740=\"],\n  [\"62\u00d781=\", \"96\u00d732=\"],\n  [\"81\u00d799=\", \"35\u00d760=\"],\n  [\"63\u00d743=\", \"46\u00d789=\"],\n  [\"56\u00d725=\", \"40\u00d717=\"],\n  [\"87\u00d793=\", \"94\u00d762=\"],\n  [\"39\u00d772=\", \"92\u00d751=\"],\n  [\"78\u00d738=\", \"69\u00d730=\"],\n  [\"39\u00d740=\", \"73\u00d792=\"],\n  [\"95\u00d735=\", \"20\u00d741=\"],\n  [\"14\u00d725=\", \"40\u00d748=\"],\n  [\"53\u00d755=\", \"56\u00d776=\"],\n  [\"97\u00d736=\", \"16\u00d762=\"],\n  [\"57\u00d763=\", \"12\u00d758=\"],\n  [\"72\u00d776=\", \"36\u00d753=\"],\n  [\"51\u00d782=\", \"52\u00d768=\"],\n  [\"85\u00d744=\", \"50\u00d799=\"],\n  [\"97\u00d788=\", \"59\u00d737=\"],\n  [\"92\u00d796=\", \"54\u00d756=\"],\n  [\"17\u00d741=\", \"95\u00d716=\"],\n  [\"82\u00d725=\", \"63\u00d733=\"],\n  [\"62\u00d725=\", \"60\u00d718=\"],\n  [\"29\u00d766=\", \"20\u00d718=\"],\n  [\"22\u00d744=\", \"78\u00d753=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the dated-worksheet update: refresh the date line and every\n# two-digit multiplication problem in the table with the new values\n# from the commit.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-04-29 Tuesday\", \"2025-04-30 Wednesday\"),\n    @(\"33\u00d774=\", \"49\u00d751=\"),\n    @(\"24\u00d742=\", \"60\u00d740=\"),\n    @(\"62\u00d781=\", \"96\u00d732=\"),\n    @(\"81\u00d799=\", \"35\u00d760=\"),\n    @(\"63\u00d743=\", \"46\u00d789=\"),\n    @(\"56\u00d725=\", \"40\u00d717=\"),\n    @(\"87\u00d793=\", \"94\u00d762=\"),\n    @(\"39\u00d772=\", \"92\u00d751=\"),\n    @(\"78\u00d738=\", \"69\u00d730=\"),\n    @(\"39\u00d740=\", \"73\u00d792=\"),\n    @(\"95\u00d735=\", \"20\u00d741=\"),\n    @(\"14\u00d725=\", \"40\u00d748=\"),\n    @(\"53\u00d755=\", \"56\u00d776=\"),\n    @(\"97\u00d736=\", \"16\u00d762=\"),\n    @(\"57\u00d763=\", \"12\u00d758=\"),\n    @(\"72\u00d776=\", \"36\u00d753=\"),\n    @(\"51\u00d782=\", \"52\u00d768=\"),\n    @(\"85\u00d744=\", \"50\u00d799=\"),\n    @(\"97\u00d788=\", \"59\u00d737=\"),\n    @(\"92\u00d796=\", \"54\u00d756=\"),\n    @(\"17\u00d741=\", \"95\u00d716=\"),\n    @(\"82\u00d725=\", \"63\u00d733=\"),\n    @(\"62\u00d725=\", \"60\u00d718=\"),\n    @(\"29\u00d766=\", \"20\u00d718=\"),\n    @(\"22\u00d744=\", \"78\u00d753=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
